# Generate Report for Handoff
#
# Moves the localization-status report from "In Translation" to
# "Ready for handoff": updates the Status columns on all three sheets,
# refreshes the "Latest Handoff"/"Latest HO Xliff Generate" timestamps,
# and widens the two "Status" columns (E/F on Overview, C on the
# per-language sheets) to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Timestamps bumped forward a few seconds ---------------------------
# Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2
# "Latest Handoff Datetime" shared the same original value.
$wsOverview.Range("G2").Value = "2016-08-26 06:57:08"
$wsDeDe.Range("H2").Value     = "2016-08-26 06:57:08"

# zh-cn!H2 "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-08-26 06:56:57"

# --- Widen the Status columns so the new text isn't truncated ---------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33   # column E
$wsOverview.Columns.Item(6).ColumnWidth = 16.33   # column F
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33   # column C
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33   # column C
